$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A4").Value = 1003
$ws.Range("B4").Value = "Test 2 Co"
$ws.Range("C4").Value = "hf"
$ws.Range("D4").Value = "PRJ-01"

$ws.Range("A5").Value = 1004
$ws.Range("B5").Value = "tessssssssssssssssssssssssssssssssssssssssssssssssssssssste lengthhhhhhhhhhhhhhhhhhhhhhhhhhhhhhhhhhhhhhhhhhhhhhhhhhhhhhhhhhhhhhh"
$ws.Range("C5").Value = 12500
$ws.Range("D5").Value = "PRJ-01"

$ws.Range("D6").Select()
